$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells are formatted as Text before assigning values so that
# values such as "1.007" or "45.90" or multi-dot prices are preserved exactly
# as strings instead of being coerced into numeric cell values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.496.81"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.840.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.75"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4631"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3845"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.90"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07878"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9946"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.45"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.857.94"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.926"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.103"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.007"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.67"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06653"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.07"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.499.95"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.371"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.87"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.069.51"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.39"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.105"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.381"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.67"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9717"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09377"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.597"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.272"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.333"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06023"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02221"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.261"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.181"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5866"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1861"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.30"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.236"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5576"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.17"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.901"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06693"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.90"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.23%  "
